# Moving calculations of measures based on single products to standalone file
#
# Inserts three new rows (21-23) of "measure" notes above the existing
# "Candidates" / "Fiddling with open tender" / "Exceptions" rows (which
# shift down from 24-26 to 27-29), fills in their content, and updates
# the sheet view selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 blank rows before the current row 21 (old rows 21-26 -> 24-29)
$ws.Rows.Item(21).Resize(3).Insert()

# 2. Row 21 - "Low-efficiency auctions"
$ws.Range("A21").Value = "Low-efficiency auctions"
$ws.Range("B21").Value = "Proportion of agency contracts awarded where price change <1%"
$ws.Range("C21").Value = "Agency"
$ws.Range("I21").Value = "Preponderance of auctions in this band difficult to justify"
$ws.Range("B21").Style = "Normal"
$ws.Cells.Item(21,2).Font.Bold = $false

# 3. Row 22 - "Nested agencies"
$ws.Range("A22").Value = "Nested agencies"
$ws.Range("B22").Value = "Is it possible to identify these measures within agences? Perhaps agency names within regnum? Or use different budget sources?"
$ws.Range("F22").Value = "My theory should apply within as across agencies. The actors are just junior and senior so we should see same patterns ""zooming in"" on agencies"
$ws.Range("I22").Value = """turtles all the way down"""

# 4. Row 23 - "Regional vs Federal"
$ws.Range("A23").Value = "Regional vs Federal"
$ws.Range("B23").Value = "Are the corruption measures different out in regions based on whether funds are federal?"
$ws.Range("F23").Value = "Identifiable from ContractFinance fields"
$ws.Range("I23").Value = "Control variable"
# (B23/F23 order confirmed against target shared-string indices 98/97)

# Apply the wrap-text style (style index 3 in the original file) to the
# descriptive columns, matching columns B, F, I used elsewhere in the sheet.
$ws.Range("B21").Style = $ws.Range("B20").Style
$ws.Range("B22").Style = $ws.Range("B20").Style
$ws.Range("B23").Style = $ws.Range("B20").Style
$ws.Range("F22").Style = $ws.Range("F17").Style
$ws.Range("F23").Style = $ws.Range("F17").Style
$ws.Range("I21").Style = $ws.Range("I20").Style
$ws.Range("I22").Style = $ws.Range("I20").Style
$ws.Range("I23").Style = $ws.Range("I20").Style

# Row heights matching the authored content (32 / 80 / 48 points)
$ws.Rows.Item(21).RowHeight = 32
$ws.Rows.Item(22).RowHeight = 80
$ws.Rows.Item(23).RowHeight = 48

# 5. Update the frozen-pane scroll position and the bottom-right selection
$aw = $excel.ActiveWindow
$aw.ScrollRow = 13
$aw.ScrollColumn = 2
$ws.Range("A23").Select()
